$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -13.376
$ws.Range("A9").Value = -21.882
$ws.Range("C9").Value = -12.42
$ws.Range("D9").Value = -7.866
$ws.Range("C11").Value = -12.126
$ws.Range("A18").Value = -21.898
$ws.Range("A20").Value = -20.636
$ws.Range("C23").Value = -12.919
$ws.Range("C24").Value = -12.38
$ws.Range("C26").Value = -12.311
$ws.Range("A27").Value = -21.859
$ws.Range("D27").Value = -7.683000000000002
$ws.Range("D29").Value = -7.709999999999999
$ws.Range("D32").Value = -7.382
$ws.Range("C34").Value = -12.049
$ws.Range("A35").Value = -20.186
$ws.Range("C35").Value = -12.12
$ws.Range("D37").Value = -7.87
$ws.Range("D38").Value = -7.873
$ws.Range("D41").Value = -7.904999999999999
$ws.Range("D45").Value = -7.525
$ws.Range("C48").Value = -11.742
$ws.Range("C49").Value = -12.996
$ws.Range("D51").Value = -8.06
$ws.Range("C52").Value = -11.95
$ws.Range("D57").Value = -8.296000000000001
$ws.Range("D64").Value = -7.887999999999998
$ws.Range("C66").Value = -11.586
$ws.Range("C67").Value = -11.59
$ws.Range("A69").Value = -21.593
$ws.Range("A76").Value = -20.241
$ws.Range("A78").Value = -20.402
$ws.Range("C78").Value = -12.822
$ws.Range("C80").Value = -12.404
$ws.Range("A82").Value = -21.801
$ws.Range("D82").Value = -7.897
$ws.Range("A83").Value = -21.84
$ws.Range("A93").Value = -21.428
$ws.Range("D93").Value = -7.214
$ws.Range("C99").Value = -11.695
$ws.Range("D102").Value = -7.811
$ws.Range("C104").Value = -12.906
$ws.Range("D105").Value = -7.876
